$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Contest 46 (row 58): MI vs DC - fill in the 9 entrant scores
$ws.Range("E58").Value = 70
$ws.Range("H58").Value = 100
$ws.Range("K58").Value = 30
$ws.Range("N58").Value = 20
$ws.Range("Q58").Value = 50
$ws.Range("T58").Value = 20
$ws.Range("W58").Value = 80
$ws.Range("Z58").Value = 60
$ws.Range("AC58").Value = 40

# Contest 47 (row 59): RR vs CSK - fill in the 9 entrant scores
$ws.Range("E59").Value = 40
$ws.Range("H59").Value = 50
$ws.Range("K59").Value = 100
$ws.Range("N59").Value = 0
$ws.Range("Q59").Value = 70
$ws.Range("T59").Value = 80
$ws.Range("W59").Value = 60
$ws.Range("Z59").Value = 20
$ws.Range("AC59").Value = 30

# In the source workbook, M58 and S58 ended up as hard-coded values rather
# than the usual VLOOKUP/RANK formula - replicate that exactly.
$ws.Range("M58").Value = -22.5
$ws.Range("S58").Value = -22.5

# Extend the season-total SUM formulas (row 71) to cover the newly-used rows
# (13:57 -> 13:68), matching the rest of the contest rows already present.
$ws.Range("E71").Formula = "=SUM(D13:D68)"
$ws.Range("H71").Formula = "=SUM(G13:G68)"
$ws.Range("K71").Formula = "=SUM(J13:J68)"
$ws.Range("N71").Formula = "=SUM(M13:M68)"
$ws.Range("Q71").Formula = "=SUM(P13:P68)"
$ws.Range("T71").Formula = "=SUM(S13:S68)"
$ws.Range("W71").Formula = "=SUM(V13:V68)"
$ws.Range("Z71").Formula = "=SUM(Y13:Y68)"
$ws.Range("AC71").Formula = "=SUM(AB13:AB68)"
